# Update QR code entrance URLs to use the new Netlify URL, and update the
# saved cell selection to C14 (matches the author's final working cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUrl = "https://uky-hospital-gps.netlify.app/start/?start="

# Update every A1:A13 formula (A2:A13 were stored as one Excel "shared
# formula" group, so re-set each row individually to keep them in sync).
for ($row = 1; $row -le 13; $row++) {
    $cell = $ws.Range("A$row")
    $cell.Formula = '=_xlfn.CONCAT("' + $newUrl + '",SUBSTITUTE(B' + $row + ',".png",""))'
}

# Update the saved selection on the sheet.
$ws.Range("C14").Select()
